$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CA-Tasks")

# Tick the four outstanding checkboxes (Deadline Met / Good Programming
# Practices Finalised / Analysis Plan confirmed / rows 6, 9, 15, 16) by
# setting their linked cells to TRUE, exactly like clicking each checkbox
# in the Excel UI would.
$ws.Range("D6").Value = $true
$ws.Range("D9").Value = $true
$ws.Range("D15").Value = $true
$ws.Range("D16").Value = $true

# Also flip the checkbox controls themselves (best effort - mirrors the
# linked-cell state for the shapes that drive D6/D9/D15/D16).
$checkboxNames = @{
    "Check Box 8"  = "D6"
    "Check Box 22" = "D9"
    "Check Box 28" = "D15"
    "Check Box 29" = "D16"
}
foreach ($name in $checkboxNames.Keys) {
    try {
        $shp = $ws.Shapes.Item($name)
        $shp.ControlFormat.Value = 1
    } catch {
    }
}

# Recalculate so the progress formula in D24 reflects the newly-checked
# items.
$wb.Application.Calculate()

# Finally, move the active selection to D3 (matches the saved cursor
# position in the workbook).
$ws.Range("D3").Select()
